$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 14.22178750299507
$ws.Range("C2").Value = 9.393696377227116
$ws.Range("D2").Value = 9.703837241471211
$ws.Range("E2").Value = 13.86929718835303
$ws.Range("F2").Value = 30.11101205632313
$ws.Range("I2").Value = 19.67133455059014
$ws.Range("J2").Value = 9.798452562305373
$ws.Range("N2").Value = 16.57561220254406
$ws.Range("O2").Value = 22.23624642048586
$ws.Range("B3").Value = 13.66403666633304
$ws.Range("C3").Value = 8.89470672416814
$ws.Range("D3").Value = 9.659539720371111
$ws.Range("E3").Value = 13.81211111198105
$ws.Range("F3").Value = 30.09216981917435
$ws.Range("I3").Value = 19.77229499921416
$ws.Range("J3").Value = 9.803993625391232
$ws.Range("N3").Value = 16.61480622921754
$ws.Range("O3").Value = 22.27421354860364
$ws.Range("B4").Value = 13.31105359881824
$ws.Range("C4").Value = 8.573939683486794
$ws.Range("D4").Value = 9.633791041157938
$ws.Range("E4").Value = 13.77976551609266
$ws.Range("F4").Value = 30.08933451730035
$ws.Range("I4").Value = 19.83900404398441
$ws.Range("J4").Value = 9.809027471025242
$ws.Range("N4").Value = 16.64065590125056
$ws.Range("O4").Value = 22.30359263973623
$ws.Range("B5").Value = 13.16477447573081
$ws.Range("C5").Value = 8.43973457570795
$ws.Range("D5").Value = 9.623670518990407
$ws.Range("E5").Value = 13.76729008282412
$ws.Range("F5").Value = 30.09037584019203
$ws.Range("I5").Value = 19.86737241629763
$ws.Range("J5").Value = 9.811489339460577
$ws.Range("N5").Value = 16.65163923458543
$ws.Range("O5").Value = 22.31708586296805
$ws.Range("B6").Value = 13.14034444435402
$ws.Range("C6").Value = 8.417243336494382
$ws.Range("D6").Value = 9.622012712074538
$ws.Range("E6").Value = 13.76526143887667
$ws.Range("F6").Value = 30.09068140754558
$ws.Range("I6").Value = 19.87215438875328
$ws.Range("J6").Value = 9.811922931155483
$ws.Range("N6").Value = 16.65349017066822
$ws.Range("O6").Value = 22.31941812057871
$ws.Range("B7").Value = 13.30909039399485
$ws.Range("C7").Value = 8.572143691080161
$ws.Range("D7").Value = 9.633653035207201
$ws.Range("E7").Value = 13.77959439841273
$ws.Range("F7").Value = 30.08933966717665
$ws.Range("I7").Value = 19.83938183980957
$ws.Range("J7").Value = 9.809059010158386
$ws.Range("N7").Value = 16.6408022058844
$ws.Range("O7").Value = 22.30376846215638
$ws.Range("B8").Value = 14.03177754272033
$ws.Range("C8").Value = 9.224699694391674
$ws.Range("D8").Value = 9.688267172098975
$ws.Range("E8").Value = 13.84901136460697
$ws.Range("F8").Value = 30.10270338333558
$ws.Range("I8").Value = 19.70516465932629
$ws.Range("J8").Value = 9.800024666501896
$ws.Range("N8").Value = 16.58875636316404
$ws.Range("O8").Value = 22.24807565219316
$ws.Range("B9").Value = 15.35742809848724
$ws.Range("C9").Value = 10.38572983715074
$ws.Range("D9").Value = 9.806518158540282
$ws.Range("E9").Value = 14.00661780036642
$ws.Range("F9").Value = 30.19810684302154
$ws.Range("I9").Value = 19.47953020030416
$ws.Range("J9").Value = 9.795239350197669
$ws.Range("N9").Value = 16.50082290419506
$ws.Range("O9").Value = 22.1871833616693
$ws.Range("B10").Value = 16.26633074300061
$ws.Range("C10").Value = 11.16169228561544
$ws.Range("D10").Value = 9.899696528500931
$ws.Range("E10").Value = 14.13481746408141
$ws.Range("F10").Value = 30.31012993904289
$ws.Range("I10").Value = 19.33682660710989
$ws.Range("J10").Value = 9.799580493783672
$ws.Range("N10").Value = 16.44479100444383
$ws.Range("O10").Value = 22.17211329713682
$ws.Range("B11").Value = 16.66413035629527
$ws.Range("C11").Value = 11.49726837044941
$ws.Range("D11").Value = 9.94333961941091
$ws.Range("E11").Value = 14.19567025463478
$ws.Range("F11").Value = 30.37010287608278
$ws.Range("I11").Value = 19.27695396906607
$ws.Range("J11").Value = 9.803253196481377
$ws.Range("N11").Value = 16.42115394043358
$ws.Range("O11").Value = 22.17173033711754
$ws.Range("B12").Value = 16.81240739323338
$ws.Range("C12").Value = 11.62179498719444
$ws.Range("D12").Value = 9.96003658024082
$ws.Range("E12").Value = 14.21906326675258
$ws.Range("F12").Value = 30.39409804200024
$ws.Range("I12").Value = 19.25501035278474
$ws.Range("J12").Value = 9.804887145941498
$ws.Range("N12").Value = 16.4124689268118
$ws.Range("O12").Value = 22.17251733069071
$ws.Range("B13").Value = 16.78057985393064
$ws.Range("C13").Value = 11.59508992327103
$ws.Range("D13").Value = 9.956433185166812
$ws.Range("E13").Value = 14.2140098697264
$ws.Range("F13").Value = 30.38887331656099
$ws.Range("I13").Value = 19.25970383035587
$ws.Range("J13").Value = 9.804524447450328
$ws.Range("N13").Value = 16.41432758491265
$ws.Range("O13").Value = 22.17230637381909
$ws.Range("B14").Value = 16.67637708329382
$ws.Range("C14").Value = 11.50756455442094
$ws.Range("D14").Value = 9.944709935969499
$ws.Range("E14").Value = 14.19758789786504
$ws.Range("F14").Value = 30.37205129698296
$ws.Range("I14").Value = 19.27513402680538
$ws.Range("J14").Value = 9.803382754978617
$ws.Range("N14").Value = 16.42043409434962
$ws.Range("O14").Value = 22.17177640003393
$ws.Range("B15").Value = 16.61223944208494
$ws.Range("C15").Value = 11.45361964621843
$ws.Range("D15").Value = 9.937550970903681
$ws.Range("E15").Value = 14.18757400959689
$ws.Range("F15").Value = 30.36191426282703
$ws.Range("I15").Value = 19.28468049168451
$ws.Range("J15").Value = 9.802715072665505
$ws.Range("N15").Value = 16.42420911185274
$ws.Range("O15").Value = 22.17157317619423
$ws.Range("B16").Value = 16.24000876108173
$ws.Range("C16").Value = 11.13940715544857
$ws.Range("D16").Value = 9.896868744271242
$ws.Range("E16").Value = 14.13089031381768
$ws.Range("F16").Value = 30.30639100576411
$ws.Range("I16").Value = 19.34084122321659
$ws.Range("J16").Value = 9.799374559633002
$ws.Range("N16").Value = 16.44637296825818
$ws.Range("O16").Value = 22.17226869389478
$ws.Range("B17").Value = 16.00756301475234
$ws.Range("C17").Value = 10.94215308315859
$ws.Range("D17").Value = 9.872225904028021
$ws.Range("E17").Value = 14.09675525868707
$ws.Range("F17").Value = 30.27463028012279
$ws.Range("I17").Value = 19.37658838957228
$ws.Range("J17").Value = 9.797759483925331
$ws.Range("N17").Value = 16.46044378699915
$ws.Range("O17").Value = 22.17435427923052
$ws.Range("B18").Value = 15.87239736474446
$ws.Range("C18").Value = 10.82706075516009
$ws.Range("D18").Value = 9.858170829934402
$ws.Range("E18").Value = 14.07736114111406
$ws.Range("F18").Value = 30.25721150532634
$ws.Range("I18").Value = 19.39762375882535
$ws.Range("J18").Value = 9.796990416609434
$ws.Range("N18").Value = 16.4687113099656
$ws.Range("O18").Value = 22.17616304516274
$ws.Range("B19").Value = 15.82638385594302
$ws.Range("C19").Value = 10.78781283687892
$ws.Range("D19").Value = 9.853432743171412
$ws.Range("E19").Value = 14.07083620188035
$ws.Range("F19").Value = 30.25145996519026
$ws.Range("I19").Value = 19.40482736995206
$ws.Range("J19").Value = 9.796757510904225
$ws.Range("N19").Value = 16.47154051335059
$ws.Range("O19").Value = 22.17688004035985
$ws.Range("B20").Value = 16.03246013238665
$ws.Range("C20").Value = 10.96332089404359
$ws.Range("D20").Value = 9.874836951209968
$ws.Range("E20").Value = 14.10036431692911
$ws.Range("F20").Value = 30.27792345541422
$ws.Range("I20").Value = 19.37273390359699
$ws.Range("J20").Value = 9.797914870988432
$ws.Range("N20").Value = 16.4589278823909
$ws.Range("O20").Value = 22.17406920429871
$ws.Range("B21").Value = 16.70704882982735
$ws.Range("C21").Value = 11.53334236010013
$ws.Range("D21").Value = 9.948148800849866
$ws.Range("E21").Value = 14.20240206964603
$ws.Range("F21").Value = 30.37695756497428
$ws.Range("I21").Value = 19.27058199476659
$ws.Range("J21").Value = 9.80371150649729
$ws.Range("N21").Value = 16.41863325390798
$ws.Range("O21").Value = 22.17190676523563
$ws.Range("B22").Value = 17.13412968299707
$ws.Range("C22").Value = 11.89101370767598
$ws.Range("D22").Value = 9.997049832952587
$ws.Range("E22").Value = 14.27111872608301
$ws.Range("F22").Value = 30.44916379081078
$ws.Range("I22").Value = 19.20806980471752
$ws.Range("J22").Value = 9.808916807395216
$ws.Range("N22").Value = 16.39384761975925
$ws.Range("O22").Value = 22.17592601907446
$ws.Range("B23").Value = 16.90748343254122
$ws.Range("C23").Value = 11.70149088636069
$ws.Range("D23").Value = 9.970863548017119
$ws.Range("E23").Value = 14.23426280502714
$ws.Range("F23").Value = 30.40994563191641
$ws.Range("I23").Value = 19.2410436802453
$ws.Range("J23").Value = 9.806009357464825
$ws.Range("N23").Value = 16.40693458666689
$ws.Range("O23").Value = 22.17328356610405
$ws.Range("B24").Value = 16.02120890112276
$ws.Range("C24").Value = 10.95375618206871
$ws.Range("D24").Value = 9.873656145821265
$ws.Range("E24").Value = 14.09873194234455
$ws.Range("F24").Value = 30.27643199067555
$ws.Range("I24").Value = 19.37447500978001
$ws.Range("J24").Value = 9.797844123719145
$ws.Range("N24").Value = 16.45961266828393
$ws.Range("O24").Value = 22.17419618745639
$ws.Range("B25").Value = 15.00963879198643
$ws.Range("C25").Value = 10.08489064325116
$ws.Range("D25").Value = 9.7733842007842
$ws.Range("E25").Value = 13.9617494639214
$ws.Range("F25").Value = 30.16490827828015
$ws.Range("I25").Value = 19.53653142203353
$ws.Range("J25").Value = 9.795151714684436
$ws.Range("N25").Value = 16.52310291909213
$ws.Range("O25").Value = 22.19845889362142

Write-Output "Updated loading_percent values for the 380 kV case (rows 2-25)."
